$d = $word.ActiveDocument

function Find-ParagraphContaining($doc, $needle) {
    foreach ($p in $doc.Paragraphs) {
        if ($p.Range.Text -like "*$needle*") {
            return $p
        }
    }
    return $null
}

# --- Step 1: the "PEN No ... : PEN" paragraph -------------------------------
# It currently holds the PEN content and uses the "spacing 276" paragraph
# formatting. Split it into:
#   (a) an emptied paragraph, centered, with no PEN content any more
#   (b) a brand-new paragraph (keeping the old "spacing 276" formatting)
#       that now carries the PEN content that used to live in (a).
$penPara = Find-ParagraphContaining $d "PEN No"

$emptiedPenXml = '<w:p><w:pPr><w:pStyle w:val="NoSpacing"/><w:jc w:val="center"/><w:rPr><w:b/><w:noProof/></w:rPr></w:pPr></w:p>'

$newPenContentXml = '<w:p><w:pPr><w:pStyle w:val="NoSpacing"/><w:spacing w:line="276" w:lineRule="auto"/><w:rPr><w:b/><w:noProof/></w:rPr></w:pPr>' +
    '<w:r><w:rPr><w:noProof/></w:rPr><w:tab/></w:r>' +
    '<w:r><w:rPr><w:b/><w:noProof/></w:rPr><w:t>PEN No</w:t></w:r>' +
    '<w:r><w:rPr><w:b/><w:noProof/></w:rPr><w:tab/></w:r>' +
    '<w:r><w:rPr><w:b/><w:noProof/></w:rPr><w:tab/><w:t xml:space="preserve">: </w:t></w:r>' +
    '<w:bookmarkStart w:id="2" w:name="PEN"/>' +
    '<w:r><w:rPr><w:b/><w:noProof/></w:rPr><w:t>PEN</w:t></w:r>' +
    '<w:bookmarkEnd w:id="2"/>' +
    '</w:p>'

$penPara.Range.InsertXML($emptiedPenXml + $newPenContentXml) | Out-Null

# --- Step 2: the "Basic Pay ... : Rs. BP/-" paragraph ----------------------
# Keep its content as-is, but append a brand-new empty paragraph right after
# it (same "spacing 276" formatting) before the table starts.
# Re-find it fresh (rather than reuse a stale pre-edit reference/index) since
# step 1 above shifted every paragraph after it down by one.
$basicPayPara = Find-ParagraphContaining $d "Basic Pay"

$basicPayXml = '<w:p><w:pPr><w:pStyle w:val="NoSpacing"/><w:spacing w:line="276" w:lineRule="auto"/><w:rPr><w:b/><w:noProof/></w:rPr></w:pPr>' +
    '<w:r><w:rPr><w:b/><w:noProof/></w:rPr><w:tab/><w:t>Basic Pay</w:t></w:r>' +
    '<w:r><w:rPr><w:b/><w:noProof/></w:rPr><w:tab/><w:t xml:space="preserve">: Rs. </w:t></w:r>' +
    '<w:bookmarkStart w:id="3" w:name="BasicPay"/>' +
    '<w:r><w:rPr><w:b/><w:noProof/></w:rPr><w:t>BP</w:t></w:r>' +
    '<w:bookmarkEnd w:id="3"/>' +
    '<w:r><w:rPr><w:b/><w:noProof/></w:rPr><w:t>/-</w:t></w:r>' +
    '</w:p>'

$newBlankXml = '<w:p><w:pPr><w:pStyle w:val="NoSpacing"/><w:spacing w:line="276" w:lineRule="auto"/><w:rPr><w:b/><w:noProof/></w:rPr></w:pPr></w:p>'

$basicPayPara.Range.InsertXML($basicPayXml + $newBlankXml) | Out-Null
